$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# Rows that get the shared "x" marker in column D
$xRows = @(2,3,5,6,8,9,10,11,13,14,16,17,19,20,22,23,25,26,28,29,31,32,34,35,37,38,44,49,52)
foreach ($r in $xRows) {
    $ws.Range("D$r").Value = "x"
}

# Rows with the new, longer tracked comments
$ws.Range("D86").Value = "Hay ticket"
$ws.Range("D87").Value = "El no poder editar no puedo cambiar el estado"
$ws.Range("D88").Value = "El no poder editar no puedo cambiar el estado"

# Update the scroll position / active selection to match the saved view state
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 79
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D90").Select()
